$wb = $excel.ActiveWorkbook

# Update sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 155
$ws1.Range("F10").Value = 5334

# Update sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 155
$ws4.Range("F10").Value = 5334
